$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select A2 and clear its contents (matches the row 2 value being removed
# while rows 3 and 4 keep their original row numbers).
$ws.Range("A2").Select()
$ws.Range("A2").ClearContents()
